$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency Price (D) / Volume(1h) (E) columns.
# A handful of the new Price strings parse as plain numbers (e.g. "1.000"),
# so those cells are pre-formatted as Text ("@") before the assignment --
# otherwise Excel would auto-convert them and drop the trailing zeros.

$ws.Range("D2").Value = '29.057.33'
$ws.Range("D3").Value = '1.826.51'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.42'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6370'
$ws.Range("E6").Value = '  -4.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.86'
$ws.Range("E8").Value = '  +6.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2938'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07348'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.82'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07673'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = '1.827.98'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.990'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6638'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.08'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008711'
$ws.Range("E17").Value = '  +5.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.048'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").Value = '29.037.85'
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '2.074.55'
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '225.78'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.40'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.135'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.000'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.71'
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.479'
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1369'
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.92'
$ws.Range("E29").Value = '  -0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.505'
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.092'
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.034'
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05322'
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.839'
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7379'
$ws.Range("E36").Value = '  -2.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.156'
$ws.Range("E37").Value = '  +2.21%  '
$ws.Range("D39").Value = '1.300.40'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01788'
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.744'
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.297'
$ws.Range("E42").Value = '  +5.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9023'
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.81'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").Value = '1.973.95'
$ws.Range("E46").Value = '  -0.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5135'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.08'
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.728'
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07248'
$ws.Range("E51").Value = '  -18.19%  '
